$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 = "Measures reported, Senate concurrent resolutions" -- fill in the
# missing House-column (C) figure.
$ws.Range("C26").Value = 1

# Row 27 = "Measures reported, House concurrent resolutions" -- fill in the
# missing House-column (C) figure.
$ws.Range("C27").Value = 5
